$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.100554585456848
$ws.Range("B1").Value = 2.111649990081787
$ws.Range("C1").Value = 9.208041191101074
$ws.Range("D1").Value = 2.409876108169556
$ws.Range("E1").Value = 1.294128775596619
